# Updates cryptos list values (Price and Volume(1h) columns) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.031.48'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '2.051.74'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '246.38'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('E6').Value = '  -1.74%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '58.74'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -4.53%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.29%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0774'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('E11').Value = '  +2.03%  '
$ws.Range('E12').Value = '  -5.39%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.891'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +7.73%  '
$ws.Range('D14').Value = '2.350.83'
$ws.Range('E14').Value = '  -0.35%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.75'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').Value = '2.078.40'
$ws.Range('E16').Value = '  +0.91%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '18.20'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '37.035.53'
$ws.Range('E18').Value = '  -0.42%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '73.96'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('D20').Value = '0.0₃0887'
$ws.Range('E20').Value = '  -2.61%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.43'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.54%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '238.09'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +0.87%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.78'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +3.74%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  -4.00%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '20.04'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.49%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '5.56'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +14.65%  '
$ws.Range('E30').Value = '  -1.12%  '
$ws.Range('E31').Value = '  -2.17%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.71'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  +6.46%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.0851'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('E39').Value = '  -1.44%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.07'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -3.20%  '
$ws.Range('E41').Value = '  -1.21%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0974'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -11.17%  '
$ws.Range('E43').Value = '  +0.47%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '97.81'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.33%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '17.00'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -6.57%  '
$ws.Range('D46').Value = '1.301.07'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('E47').Value = '  -4.29%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.87'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -0.60%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '6.75'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('D50').Value = '2.234.47'
$ws.Range('E50').Value = '  -0.55%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.59'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.25%  '
